$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: fill in the first transaction line (previously blank placeholder row) ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ORGASOL LIGHT CREAM"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "'0"
$ws.Range("N7").Value = "130.00"
$ws.Range("P7").Value = "130.0000"
$ws.Range("Q7").Value = "'1:0"

# --- Insert a new row at 8 for the second transaction line, pushing the totals/footer down ---
$ws.Rows.Item(8).Insert()

# Row 8 merges mirror row 7's layout
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "PRISBRINA  CAPS"
$ws.Range("H8").Value = "-1:-1"
$ws.Range("L8").Value = "'0"
$ws.Range("N8").Value = "150.00"
$ws.Range("P8").Value = "150.0000"
$ws.Range("Q8").Value = "'1:0"

# Match row7 styling for the newly inserted row8 cells
$ws.Range("A8:B8").Style = $ws.Range("A7:B7").Style
$ws.Range("C8:G8").Style = $ws.Range("C7:G7").Style
$ws.Range("H8:K8").Style = $ws.Range("H7:K7").Style
$ws.Range("L8:M8").Style = $ws.Range("L7:M7").Style
$ws.Range("N8:O8").Style = $ws.Range("N7:O7").Style
$ws.Range("P8").Style = $ws.Range("P7").Style
$ws.Range("Q8").Style = $ws.Range("Q7").Style
$ws.Rows.Item(8).RowHeight = 24.75

# --- Row 9 (former row 8 with the P/Q total cells) now holds the transactions-count total ---
$ws.Range("P9").Value = 280

# --- Row 10 (former row 9, footer) gets the updated timestamp ---
$ws.Range("A10").Value = "Saturday, 24 May, 2025 9:42 AM"

Write-Host "edit applied"
